# Update column F (dSF) values on Sheet1 as part of a data repull / mean calculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = 2
    3  = -1
    6  = 5
    7  = 6
    8  = 7
    9  = 5
    10 = -1
    11 = 2
    12 = 5
    13 = 3
    14 = 3
    15 = 2
    16 = 6
    17 = -2
    18 = 3
    19 = 6
    20 = 3
    21 = 4
    22 = -2
    23 = 1
    24 = -2
    25 = -1
    26 = 3
    27 = -5
    28 = -1
    29 = 2
    30 = -1
    31 = 1
    33 = -5
    34 = -2
    35 = 4
    36 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
